# Scheduled-runner style update: refresh cached market-price / profit
# figures (columns H-N) on several leve rows across multiple job sheets.
# Values are plain cached numbers (no formulas in this workbook), so each
# changed cell is written directly. A couple of rows drop a stale
# LeveProfit cell entirely (ClearContents), matching the source data feed.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 878.3913
$ws.Range("I28").Value = 752.375
$ws.Range("J28").Value = 1166.4286
$ws.Range("K28").Value = 752.375
$ws.Range("L28").Value = 1166.4286
$ws.Range("M28").Value = -267.375
$ws.Range("N28").Value = -2136.4286

$ws.Range("H33").Value = 277.41666
$ws.Range("I33").Value = 215.44444
$ws.Range("J33").Value = 463.33334
$ws.Range("K33").Value = 215.44444
$ws.Range("L33").Value = 463.33334
$ws.Range("M33").Value = 13.55556000000001
$ws.Range("N33").Value = -921.33334

$ws.Range("H62").Value = 4921.467
$ws.Range("I62").Value = 4509.8184
$ws.Range("J62").Value = 6053.5
$ws.Range("K62").Value = 4509.8184
$ws.Range("L62").Value = 6053.5
$ws.Range("M62").Value = -3885.8184
$ws.Range("N62").Value = -7301.5

$ws.Range("H65").Value = 4921.467
$ws.Range("I65").Value = 4509.8184
$ws.Range("J65").Value = 6053.5
$ws.Range("K65").Value = 22549.092
$ws.Range("L65").Value = 30267.5
$ws.Range("M65").Value = -19429.092
$ws.Range("N65").Value = -36507.5

$ws.Range("H116").Value = 783867
$ws.Range("J116").Value = 5309.5
$ws.Range("L116").Value = 5309.5
$ws.Range("N116").Value = -12193.5

$ws.Range("H137").Value = 21253.312
$ws.Range("I137").Value = 14667.458
$ws.Range("J137").Value = 41010.875
$ws.Range("K137").Value = 44002.374
$ws.Range("L137").Value = 123032.625
$ws.Range("M137").Value = -41452.374
$ws.Range("N137").Value = -128132.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13580.813
$ws.Range("I32").Value = 13761.071
$ws.Range("J32").Value = 6010
$ws.Range("K32").Value = 13761.071
$ws.Range("L32").Value = 6010
$ws.Range("M32").Value = -13474.071
$ws.Range("N32").Value = -6584

$ws.Range("H45").Value = 4912.0713
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 4912.0713
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 4912.0713
$ws.Range("N45").Value = -5666.0713
$ws.Range("M45").ClearContents()

$ws.Range("H74").Value = 143701.14
$ws.Range("I74").Value = 188729.03
$ws.Range("J74").Value = 12710.909
$ws.Range("K74").Value = 188729.03
$ws.Range("L74").Value = 12710.909
$ws.Range("M74").Value = -187855.03
$ws.Range("N74").Value = -14458.909

$ws.Range("H77").Value = 143701.14
$ws.Range("I77").Value = 188729.03
$ws.Range("J77").Value = 12710.909
$ws.Range("K77").Value = 943645.15
$ws.Range("L77").Value = 63554.545
$ws.Range("M77").Value = -939277.15
$ws.Range("N77").Value = -72290.545

$ws.Range("H102").Value = 1604.6471
$ws.Range("I102").Value = 845.2
$ws.Range("K102").Value = 845.2
$ws.Range("M102").Value = 776.8

$ws.Range("H132").Value = 1835.6945
$ws.Range("I132").Value = 1440.1562
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 4320.4686
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -1790.4686
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1690.5
$ws.Range("I20").Value = 1291.6154
$ws.Range("J20").Value = 2431.2856
$ws.Range("K20").Value = 1291.6154
$ws.Range("L20").Value = 2431.2856
$ws.Range("M20").Value = -1044.6154
$ws.Range("N20").Value = -2925.2856

$ws.Range("H105").Value = 2023.7667
$ws.Range("I105").Value = 1827
$ws.Range("K105").Value = 1827
$ws.Range("M105").Value = -80

$ws.Range("H107").Value = 2721.1738
$ws.Range("I107").Value = 2721.1738
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2721.1738
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -801.1738
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2159.625
$ws.Range("I16").Value = 1700
$ws.Range("J16").Value = 2619.25
$ws.Range("K16").Value = 1700
$ws.Range("L16").Value = 2619.25
$ws.Range("M16").Value = -1413
$ws.Range("N16").Value = -3193.25

$ws.Range("H31").Value = 2176791.8
$ws.Range("I31").Value = 3228543
$ws.Range("J31").Value = 3172.1333
$ws.Range("K31").Value = 3228543
$ws.Range("L31").Value = 3172.1333
$ws.Range("M31").Value = -3228248
$ws.Range("N31").Value = -3762.1333

$ws.Range("H34").Value = 2176791.8
$ws.Range("I34").Value = 3228543
$ws.Range("J34").Value = 3172.1333
$ws.Range("K34").Value = 3228543
$ws.Range("L34").Value = 3172.1333
$ws.Range("M34").Value = -3228341
$ws.Range("N34").Value = -3576.1333

$ws.Range("H113").Value = 2159.625
$ws.Range("I113").Value = 1700
$ws.Range("J113").Value = 2619.25
$ws.Range("K113").Value = 1700
$ws.Range("L113").Value = 2619.25
$ws.Range("M113").Value = 470
$ws.Range("N113").Value = -6959.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 4154.933
$ws.Range("I68").Value = 1933
$ws.Range("J68").Value = 4599.32
$ws.Range("K68").Value = 5799
$ws.Range("L68").Value = 13797.96
$ws.Range("M68").Value = -4988
$ws.Range("N68").Value = -15419.96

$ws.Range("H71").Value = 4154.933
$ws.Range("I71").Value = 1933
$ws.Range("J71").Value = 4599.32
$ws.Range("K71").Value = 17397
$ws.Range("L71").Value = 41393.88
$ws.Range("M71").Value = -13341
$ws.Range("N71").Value = -49505.88

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 47362.363
$ws.Range("I102").Value = 57254.555
$ws.Range("J102").Value = 2847.5
$ws.Range("K102").Value = 57254.555
$ws.Range("L102").Value = 2847.5
$ws.Range("M102").Value = -55632.555
$ws.Range("N102").Value = -6091.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1564.4445
$ws.Range("I61").Value = 1290.7333
$ws.Range("J61").Value = 2933
$ws.Range("K61").Value = 1290.7333
$ws.Range("L61").Value = 2933
$ws.Range("M61").Value = -1088.7333
$ws.Range("N61").Value = -3337

$ws.Range("H113").Value = 1564.4445
$ws.Range("I113").Value = 1290.7333
$ws.Range("J113").Value = 2933
$ws.Range("K113").Value = 1290.7333
$ws.Range("L113").Value = 2933
$ws.Range("M113").Value = 879.2666999999999
$ws.Range("N113").Value = -7273

$ws.Range("H136").Value = 9624.375
$ws.Range("I136").Value = 7832.6665
$ws.Range("J136").Value = 14999.5
$ws.Range("K136").Value = 23497.9995
$ws.Range("L136").Value = 44998.5
$ws.Range("M136").Value = -20947.9995
$ws.Range("N136").Value = -50098.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1057.6923
$ws.Range("I107").Value = 1012.2381
$ws.Range("J107").Value = 1248.6
$ws.Range("K107").Value = 3036.7143
$ws.Range("L107").Value = 3745.8
$ws.Range("M107").Value = -1116.7143
$ws.Range("N107").Value = -7585.799999999999
